# Add a "solidity" column to the "geometry" sheet, between the existing
# "pitch_chord_ratio" column (AM) and "thickness_max_chord_ratio" column
# (formerly AN, now shifted to AO).
#
# Inserting a whole column shifts every column at/after it one position to
# the right (and Excel carries the header cell's style along), which is
# exactly what's needed to turn A1:AQ2 into A1:AR2 while leaving all other
# existing columns' content untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("geometry")

# Column 40 is "AN" - the first column whose header changes position in the diff.
$ws.Columns.Item(40).Insert()

$ws.Cells.Item(1, 40).Value = "solidity"
$ws.Cells.Item(2, 40).Value = "[1.42997704 1.70997375]"
